$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet used to have two header rows (row1 + row2). The new layout
# consolidates them into a single header row. Deleting the old row 2
# shifts all the data rows (old 3..16) up by one (new 2..15) and keeps
# their values/styles intact.
$ws.Rows.Item(2).Delete()

# Rewrite the (now single) header row with the new column headers.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 should carry no special formatting (default style).
$ws.Range("A1:E1").Style = "Normal"

# F1:K1 use a dedicated header style (same font as the rest of the
# table, general number format).
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").WrapText = $false

# Update the sheet view: the selection now targets the header row
# (A2:K2) instead of the stale K23 reference.
$ws.Range("A2:K2").Select() | Out-Null
